$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 2250
$ws.Range("I8").Value = 250.33333
$ws.Range("K8").Value = 750.99999
$ws.Range("M8").Value = -611.99999
# Row 17
$ws.Range("H17").Value = 666.5
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2832
# Row 42
$ws.Range("H42").Value = 247.66667
$ws.Range("I42").Value = 174.66667
$ws.Range("J42").Value = 466.66666
$ws.Range("K42").Value = 524.00001
$ws.Range("L42").Value = 1399.99998
$ws.Range("M42").Value = -294.00001
$ws.Range("N42").Value = -1859.99998
# Row 93
$ws.Range("H93").Value = 55000
$ws.Range("J93").Value = 55000
$ws.Range("L93").Value = 55000
$ws.Range("N93").Value = -59992
# Row 115
$ws.Range("H115").Value = 1367.5
$ws.Range("I115").Value = 1367.5
$ws.Range("K115").Value = 4102.5
$ws.Range("M115").Value = -2535.5
# Row 127
$ws.Range("H127").Value = 1537.2
$ws.Range("I127").Value = 1757.6666
$ws.Range("K127").Value = 5272.9998
$ws.Range("M127").Value = -312.9997999999996
# Row 132
$ws.Range("H132").Value = 5338
$ws.Range("I132").Value = 3299.25
$ws.Range("K132").Value = 9897.75
$ws.Range("M132").Value = -7367.75
# Row 137
$ws.Range("H137").Value = 721357.0600000001
$ws.Range("I137").Value = 1254062.4
$ws.Range("J137").Value = 11083.333
$ws.Range("K137").Value = 3762187.2
$ws.Range("L137").Value = 33249.999
$ws.Range("M137").Value = -3759637.2
$ws.Range("N137").Value = -38349.999
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1749.4286
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 425
$ws.Range("J22").Value = 550
$ws.Range("L22").Value = 550
$ws.Range("N22").Value = -896

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
# Row 22
$ws.Range("H22").Value = 498.6
$ws.Range("J22").Value = 847
$ws.Range("L22").Value = 847
$ws.Range("N22").Value = -1547
# Row 58
$ws.Range("H58").Value = 8228
$ws.Range("I58").Value = 1399.25
$ws.Range("K58").Value = 1399.25
$ws.Range("M58").Value = -1196.25
# Row 74
$ws.Range("H74").Value = 24875
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 24750
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 24750
$ws.Range("M74").Value = -24126
$ws.Range("N74").Value = -26498
# Row 77
$ws.Range("H77").Value = 24875
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 24750
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 74250
$ws.Range("M77").Value = -70632
$ws.Range("N77").Value = -82986
# Row 99
$ws.Range("H99").Value = 6669
$ws.Range("I99").Value = 9006
$ws.Range("K99").Value = 9006
$ws.Range("M99").Value = -7508
# Row 126
$ws.Range("H126").Value = 6669
$ws.Range("I126").Value = 9006
$ws.Range("K126").Value = 27018
$ws.Range("M126").Value = -24548
# Row 132
$ws.Range("H132").Value = 15499.833
$ws.Range("J132").Value = 17999.8
$ws.Range("L132").Value = 53999.39999999999
$ws.Range("N132").Value = -59059.39999999999
# Row 134
$ws.Range("H134").Value = 12574.75
# Row 136
$ws.Range("H136").Value = 8228
$ws.Range("I136").Value = 1399.25
$ws.Range("K136").Value = 4197.75
$ws.Range("M136").Value = -1647.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 18415
$ws.Range("I3").Value = 18415
$ws.Range("K3").Value = 55245
$ws.Range("M3").Value = -55133
# Row 4
$ws.Range("H4").Value = 35521296
$ws.Range("J4").Value = 64382130
$ws.Range("L4").Value = 193146390
$ws.Range("N4").Value = -193146614
# Row 26
$ws.Range("H26").Value = 13499
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 13499
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 40497
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -41073
# Row 68
$ws.Range("H68").Value = 1100
$ws.Range("I68").Value = 650
$ws.Range("K68").Value = 1950
$ws.Range("M68").Value = -1139
# Row 71
$ws.Range("H71").Value = 1100
$ws.Range("I71").Value = 650
$ws.Range("K71").Value = 5850
$ws.Range("M71").Value = -1794
# Row 87
$ws.Range("H87").Value = 200
$ws.Range("I87").Value = 200
$ws.Range("K87").Value = 600
$ws.Range("M87").Value = 648
# Row 90
$ws.Range("H90").Value = 200
$ws.Range("I90").Value = 200
$ws.Range("K90").Value = 1800
$ws.Range("M90").Value = 4440
# Row 118
$ws.Range("H118").Value = 864
$ws.Range("I118").Value = 864
$ws.Range("K118").Value = 2592
$ws.Range("M118").Value = -1349
# Row 140
$ws.Range("H140").Value = 2024.5714
$ws.Range("I140").Value = 1862
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 5586
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -406
$ws.Range("N140").Value = -19360
# Row 141
$ws.Range("H141").Value = 230
$ws.Range("I141").Value = 230
$ws.Range("K141").Value = 690
$ws.Range("M141").Value = 4490

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 74.833336
$ws.Range("J2").Value = 100.5
$ws.Range("L2").Value = 100.5
$ws.Range("N2").Value = -326.5
# Row 11
$ws.Range("H11").Value = 2999999.2
$ws.Range("I11").Value = 2999999.2
$ws.Range("K11").Value = 2999999.2
$ws.Range("M11").Value = -2999860.2
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 132
$ws.Range("H132").Value = 125992
$ws.Range("I132").Value = 139857.6
$ws.Range("K132").Value = 419572.8
$ws.Range("M132").Value = -417042.8
# Row 141
$ws.Range("H141").Value = 79201.42999999999
$ws.Range("J141").Value = 79201.42999999999
$ws.Range("L141").Value = 79201.42999999999
$ws.Range("N141").Value = -89561.42999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 29999.723
$ws.Range("J2").Value = 99998.75
$ws.Range("L2").Value = 99998.75
$ws.Range("N2").Value = -100222.75
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 100
$ws.Range("H100").Value = 1100
$ws.Range("I100").Value = 1100
$ws.Range("K100").Value = 1100
$ws.Range("M100").Value = -559
# Row 132
$ws.Range("H132").Value = 16833.166
$ws.Range("J132").Value = 16833.166
$ws.Range("L132").Value = 50499.49800000001
$ws.Range("N132").Value = -55559.49800000001
# Row 136
$ws.Range("H136").Value = 10049.6
$ws.Range("I136").Value = 2499.8333
$ws.Range("J136").Value = 21374.25
$ws.Range("K136").Value = 7499.499899999999
$ws.Range("L136").Value = 64122.75
$ws.Range("M136").Value = -4949.499899999999
$ws.Range("N136").Value = -69222.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 9000
$ws.Range("I32").Value = 9000
$ws.Range("K32").Value = 9000
$ws.Range("M32").Value = -8683
# Row 62
$ws.Range("H62").Value = 2250
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3748
# Row 65
$ws.Range("H65").Value = 2250
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -18740
# Row 125
$ws.Range("H125").Value = 75997
$ws.Range("J125").Value = 75997
$ws.Range("L125").Value = 75997
$ws.Range("N125").Value = -85837
# Row 141
$ws.Range("H141").Value = 99995.664
$ws.Range("J141").Value = 99995.664
$ws.Range("L141").Value = 99995.664
$ws.Range("N141").Value = -110355.664
